$wb = $excel.ActiveWorkbook

$retailer  = $wb.Worksheets.Item("Retailer")
$vehicleE1 = $wb.Worksheets.Item("VehicleE1")

# Insert a new worksheet ("Sheet1") right before VehicleE1 - this is where
# the retailers that exceed the vehicle-count constraint ("m") are moved to.
$newSheet = $wb.Worksheets.Add($vehicleE1)

# Move rows 12-21 (retailers R11..R20) from Retailer into the new sheet,
# starting at A1 (no header row there).
$srcRange = $retailer.Range("A12:G21")
$srcRange.Copy()
$destRange = $newSheet.Range("A1")
$newSheet.Paste($destRange)

# Re-apply the number format used on column G (it's lost on plain Paste).
$newSheet.Range("G1:G10").NumberFormat = $retailer.Range("G12").NumberFormat

# Remove the now-duplicated rows from Retailer; this also shrinks its
# dimension/used-range back down to A1:G11.
$retailer.Range("A12:G21").EntireRow.Delete()

# Restore/update the selections seen in each sheet.
$retailer.Activate()
$null = $retailer.Rows("12:16").Select()

$newSheet.Activate()
$null = $newSheet.Range("A1:XFD5").Select()
